$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 101 (pushes old rows 101-108 down to 102-109),
# matching new tests being added above the PingPong block.
$ws.Rows.Item(101).Insert()

# New row 101: Max_Instances_4 / Feature Integration ???
$ws.Range("A101").Value = "Max_Instances_4"
$ws.Range("B101").Value = "Feature Integration ???"
# The row-insert copies formatting (and stray empty cells) from the row
# above into the rest of the new row; clear that so only A101/B101 exist.
$ws.Range("C101:H101").Clear()

# Three new rows appended at the bottom of the sheet (110-112):
# two-phase-commit_1, Multi_Paxos_3, Multi_Paxos_4
$ws.Range("A110").Value = "two-phase-commit_1"
$ws.Range("B110").Value = "Feature Integration: protocol sample"
$ws.Range("C110").Value = "No "
$ws.Range("D110").Value = "No"
$ws.Range("E110").Value = "Yes"

$ws.Range("A111").Value = "Multi_Paxos_3"
$ws.Range("B111").Value = "Feature Integration: protocol sample"
$ws.Range("C111").Value = "No "
$ws.Range("D111").Value = "Yes"

$ws.Range("A112").Value = "Multi_Paxos_4"
$ws.Range("B112").Value = "Feature Integration: protocol sample"
$ws.Range("C112").Value = "No "
$ws.Range("D112").Value = "Yes"

# "Also added as syntactic check only" note, written last so it becomes
# the final new shared-string entry.
$ws.Range("H110").Value = "Also added as syntactic check only"
$ws.Range("H111").Value = "Also added as syntactic check only"
$ws.Range("H112").Value = "Also added as syntactic check only"

# Update the sheet view's active cell / selection to match, and scroll
# the view down towards the newly added rows.
$ws.Range("H112").Select()
$excel.ActiveWindow.ScrollRow = 86
$excel.ActiveWindow.ScrollColumn = 1
